$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data: artlib, artprix, artdispo, catlib, artdesc
$rows = @(
    @("La Joconde", 410.0, 4, "classique", "Lisa, aussi connue sous le nom de Mona Lisa, Lisa di Antonio Maria Gherardini et de Lisa del Giocondo en italien, est un membre de la famille Gherardini de Florence"),
    @("Gustav Klimt", 466.0, 50, "personnes", "Gustav Klimt, né le 14 juillet 1862 à Baumgarten en Autriche et mort le 6 février 1918 à Vienne."),
    @("La Nuit étoilée", 200.0, 54, "paysages", "La Nuit étoilée est une peinture de l'artiste peintre postimpressionniste néerlandais Vincent van Gogh."),
    @("test", 120.0, 15, "paysages", "regeprgoerrtre")
)

$startRow = 3
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}
